$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Extend header merge A1:H1 -> A1:J1 (columns I and J added to the table)
# ---------------------------------------------------------------------------
$ws.Range("B1:C1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("A1:J1").Merge()

# ---------------------------------------------------------------------------
# 2. Row 4 header: two new artifact-point columns (6 and 7), and the
#    "Esfuerzo porcentaje" / "Nota individual" headers move from G4/H4
#    to I4/J4.
# ---------------------------------------------------------------------------
$ws.Range("G4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("H4").Copy()
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("I4").Value = "Esfuerzo porcentaje"
$ws.Range("J4").Value = "Nota individual"

$ws.Range("B4").Copy()
$ws.Range("G4").PasteSpecial(-4122)
$ws.Range("H4").PasteSpecial(-4122)
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 7

# ---------------------------------------------------------------------------
# 3. Member contribution grid (rows 5-8) re-marked against the new 7-point
#    scale, plus the percentage/score columns shifted from G:H to I:J.
# ---------------------------------------------------------------------------

# -- Row 5: Vivian Gomez --
$ws.Range("B4").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("H5").Copy()
$ws.Range("J5").PasteSpecial(-4122)

$ws.Range("B5").Value = "x"
$ws.Range("C5").Value = "x"
$ws.Range("D5").Value = "x"
$ws.Range("E5").Value = "x"
$ws.Range("F5").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("H5").ClearContents()
$ws.Range("F5").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("I5").Value = 34
$ws.Range("J5").Value = 5

# -- Row 6: Michael Osorio --
$ws.Range("B4").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("H5").Copy()
$ws.Range("J6").PasteSpecial(-4122)

$ws.Range("B6").Value = "x"
$ws.Range("D6").ClearContents()
$ws.Range("E6").Value = "x"
$ws.Range("F6").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("F6").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("H6").Value = "x"
$ws.Range("I6").Value = 22
$ws.Range("J6").Value = 5

# -- Row 7: Esteban Reyes --
$ws.Range("B4").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("H5").Copy()
$ws.Range("J7").PasteSpecial(-4122)

$ws.Range("B7").Value = "x"
$ws.Range("D7").ClearContents()
$ws.Range("E7").Value = "x"
$ws.Range("H7").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").Value = 22
$ws.Range("J7").Value = 5

# -- Row 8: Alvaro Yepes --
$ws.Range("B4").Copy()
$ws.Range("I8").PasteSpecial(-4122)
$ws.Range("H5").Copy()
$ws.Range("J8").PasteSpecial(-4122)

$ws.Range("B8").Value = "x"
$ws.Range("C8").ClearContents()
$ws.Range("E8").Value = "x"
$ws.Range("H8").Copy()
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("G8").Value = "x"
$ws.Range("H8").ClearContents()
$ws.Range("I8").Value = 22
$ws.Range("J8").Value = 5

# -- Row 9: TOTAL --
$ws.Range("C9").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$ws.Range("J9").PasteSpecial(-4122)
$ws.Range("J9").Value = 5

# ---------------------------------------------------------------------------
# 4. Blank formatting-only fill for the new I/J columns on rows that don't
#    carry a value (rows 2-3 and the artefact table rows 10-15), matching
#    the pre-existing style of the corresponding columns in those rows.
# ---------------------------------------------------------------------------
$ws.Range("A2:B2").Copy()
$ws.Range("I2:J2").PasteSpecial(-4122)

$ws.Range("C3:D3").Copy()
$ws.Range("I3:J3").PasteSpecial(-4122)

$ws.Range("A10:B10").Copy()
$ws.Range("I10:J10").PasteSpecial(-4122)

$ws.Range("G11:H11").Copy()
$ws.Range("I11:J11").PasteSpecial(-4122)

$ws.Range("G12:H12").Copy()
$ws.Range("I12:J12").PasteSpecial(-4122)

$ws.Range("G13:H13").Copy()
$ws.Range("I13:J13").PasteSpecial(-4122)

$ws.Range("G14:H14").Copy()
$ws.Range("I14:J14").PasteSpecial(-4122)

$ws.Range("G15:H15").Copy()
$ws.Range("I15:J15").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5. Artefacto list (rows 12-16) renamed, two new artefacts appended as
#    rows 17-18 (points 6 and 7).
# ---------------------------------------------------------------------------
$ws.Range("B12").Value = "Diseno front"
$ws.Range("B13").Value = "SAD"
$ws.Range("B14").Value = "poster"
$ws.Range("B15").Value = "Nuevos servicios back"
$ws.Range("B16").Value = "Video"

$ws.Range("A16:B16").Copy()
$ws.Range("A17:B17").PasteSpecial(-4122)
$ws.Range("A18:B18").PasteSpecial(-4122)
$ws.Range("A17").EntireRow.RowHeight = 15.75
$ws.Range("A18").EntireRow.RowHeight = 15.75

$ws.Range("A17").Value = 6
$ws.Range("B17").Value = "Diagrama de desarrollo"
$ws.Range("A18").Value = 7
$ws.Range("B18").Value = "Arreglos integraci[on back-front"

# ---------------------------------------------------------------------------
# 6. Selection, matching the saved cursor position in the target file.
# ---------------------------------------------------------------------------
$ws.Range("E11").Select()
